# Generate Report for Handoff
#
# The three "in progress" / "ready for handoff" items (34d0a952, cc596870,
# e55d1b7b) are re-sorted across rows 7-9 of every sheet (Overview, zh-cn,
# de-de). e55d1b7b finished its handback/handoff cycle, so its status moves
# from "In Translation" to "Ready for handoff" along with an updated
# timestamp, and the three rows are re-ordered: row7 <- 34d0a952 data,
# row8 <- cc596870 data, row9 <- e55d1b7b data (with refreshed datetimes).

$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        row = 7
        md = "34d0a952-eade-4acf-8ef8-b8771c6a7fdb.md"
        path = "e2e\34d0a952-eade-4acf-8ef8-b8771c6a7fdb.md"
        status = "Ready for handoff"
        ovDate = "2016-08-25 10:43:32"
        zhXlf = "34d0a952-eade-4acf-8ef8-b8771c6a7fdb.ae1b18336d55ec83182504dc2e858c5ee06988a6.zh-cn.xlf"
        zhDate = "2016-08-25 10:43:28"
        deXlf = "34d0a952-eade-4acf-8ef8-b8771c6a7fdb.ae1b18336d55ec83182504dc2e858c5ee06988a6.de-de.xlf"
        deDate = "2016-08-25 10:43:32"
    },
    @{
        row = 8
        md = "cc596870-40a6-47a8-b633-3899f55281db.md"
        path = "e2e\cc596870-40a6-47a8-b633-3899f55281db.md"
        status = "Ready for handoff"
        ovDate = "2016-08-25 10:41:06"
        zhXlf = "cc596870-40a6-47a8-b633-3899f55281db.45040b68c4809a16138982fa85dc5157999c66a7.zh-cn.xlf"
        zhDate = "2016-08-25 10:40:56"
        deXlf = "cc596870-40a6-47a8-b633-3899f55281db.45040b68c4809a16138982fa85dc5157999c66a7.de-de.xlf"
        deDate = "2016-08-25 10:41:06"
    },
    @{
        row = 9
        md = "e55d1b7b-8b0f-4bfd-817e-160c1b906dc5.md"
        path = "e2e\e55d1b7b-8b0f-4bfd-817e-160c1b906dc5.md"
        status = "Ready for handoff"
        ovDate = "2016-08-25 10:46:45"
        zhXlf = "e55d1b7b-8b0f-4bfd-817e-160c1b906dc5.87e2ac9bbb8e578c3b5c79e3bb9c6f0e733f4649.zh-cn.xlf"
        zhDate = "2016-08-25 10:46:41"
        deXlf = "e55d1b7b-8b0f-4bfd-817e-160c1b906dc5.87e2ac9bbb8e578c3b5c79e3bb9c6f0e733f4649.de-de.xlf"
        deDate = "2016-08-25 10:46:45"
    }
)

# ---- Overview sheet: columns A (file), B (path/hyperlink), E/F (status), G (date) ----
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("A" + $r.row).Value = $r.md
    $wsOverview.Range("B" + $r.row).Value = $r.path
    $wsOverview.Range("E" + $r.row).Value = $r.status
    $wsOverview.Range("F" + $r.row).Value = $r.status
    $wsOverview.Range("G" + $r.row).Value = $r.ovDate
}
foreach ($hl in $wsOverview.Hyperlinks) {
    foreach ($r in $rows) {
        if ($hl.Range.Row -eq $r.row -and $hl.Range.Column -eq 2) {
            $hl.TextToDisplay = $r.path
        }
    }
}

# ---- zh-cn sheet: columns A (file), C (status), G (target xlf), H (target date) ----
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("A" + $r.row).Value = $r.md
    $wsZh.Range("C" + $r.row).Value = $r.status
    $wsZh.Range("G" + $r.row).Value = $r.zhXlf
    $wsZh.Range("H" + $r.row).Value = $r.zhDate
}
foreach ($hl in $wsZh.Hyperlinks) {
    foreach ($r in $rows) {
        if ($hl.Range.Row -eq $r.row -and $hl.Range.Column -eq 1) {
            $hl.TextToDisplay = $r.md
        }
    }
}

# ---- de-de sheet: columns A (file), C (status), G (target xlf), H (target date) ----
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("A" + $r.row).Value = $r.md
    $wsDe.Range("C" + $r.row).Value = $r.status
    $wsDe.Range("G" + $r.row).Value = $r.deXlf
    $wsDe.Range("H" + $r.row).Value = $r.deDate
}
foreach ($hl in $wsDe.Hyperlinks) {
    foreach ($r in $rows) {
        if ($hl.Range.Row -eq $r.row -and $hl.Range.Column -eq 1) {
            $hl.TextToDisplay = $r.md
        }
    }
}
